$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2025-12-01 18:39:07"

for ($r = 2; $r -le 23; $r++) {
    $ws.Cells.Item($r, 1).Value = $newTimestamp
}
